$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "CustomerType" column (G) with header + values for existing + new rows
$ws.Range("G1").Value = "CustomerType"
$ws.Range("G2").Value = "Individual"

# New customer row: "Ababio and Sons" (a Company)
$ws.Range("A3").Value = "Ababio and Sons"
$ws.Range("B3").Value = "ababioandsons@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:ababioandsons@gmail.com") | Out-Null
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("C3").Value = 542542299
$ws.Range("E3").Value = "Accra"
$ws.Range("F3").Value = "GD-898-0909"
$ws.Range("G3").Value = "Company"

# Column G width, matching the author's template formatting
# (ColumnWidth is in "Normal"-font characters; Excel internally pads/rounds
#  this to the stored <col width> grid units, so 20.1 -> stored width 21)
$ws.Columns.Item(7).ColumnWidth = 20.1

# Selection left where the author's cursor ended up after editing
$ws.Range("G18").Select()
